# Apply the "cryptos list" data refresh (GitHub Actions scheduled update).
# Only columns B (Coin), C (Link), D (Price) and E (Volume/1h) change; column A
# (rank index) is untouched. Several Price values look numeric (e.g. "0.998")
# but must stay stored as text, exactly like the rest of column D, so they are
# assigned with a leading apostrophe (quote-prefix) which Excel strips while
# keeping the cell as Text/General instead of auto-converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.128.48'
$ws.Range('E2').Value = '  -1.68%  '

$ws.Range('E3').Value = '  -2.97%  '

$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.41%  '

$ws.Range('D5').Value = '''227.43'
$ws.Range('E5').Value = '  -2.58%  '

$ws.Range('D6').Value = '''0.608'
$ws.Range('E6').Value = '  -4.30%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '''55.17'
$ws.Range('E8').Value = '  -5.14%  '

$ws.Range('D9').Value = '''0.382'
$ws.Range('E9').Value = '  -2.74%  '

$ws.Range('D10').Value = '''0.0792'
$ws.Range('E10').Value = '  +1.43%  '

$ws.Range('E11').Value = '  -3.63%  '

$ws.Range('D12').Value = '2.323.62'
$ws.Range('E12').Value = '  -2.92%  '

$ws.Range('D13').Value = '''14.31'
$ws.Range('E13').Value = '  -5.71%  '

$ws.Range('D14').Value = '''20.48'
$ws.Range('E14').Value = '  -2.98%  '

$ws.Range('D15').Value = '''0.745'
$ws.Range('E15').Value = '  -4.39%  '

$ws.Range('E16').Value = '  -3.28%  '

$ws.Range('D17').Value = '2.017.26'
$ws.Range('E17').Value = '  -3.02%  '

$ws.Range('D18').Value = '37.020.22'
$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('E19').Value = '  -1.27%  '

$ws.Range('D20').Value = '''68.78'
$ws.Range('E20').Value = '  -3.07%  '

$ws.Range('D21').Value = '0.0₃0836'
$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('D22').Value = '''223.08'
$ws.Range('E22').Value = '  -2.75%  '

$ws.Range('E23').Value = '  +0.27%  '

$ws.Range('D24').Value = '''2.38'
$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('E25').Value = '  -5.36%  '

$ws.Range('D26').Value = '''9.35'
$ws.Range('E26').Value = '  -3.80%  '

$ws.Range('D27').Value = '''167.64'
$ws.Range('E27').Value = '  -1.89%  '

$ws.Range('E28').Value = '  -6.53%  '

$ws.Range('D29').Value = '''18.69'
$ws.Range('E29').Value = '  -4.27%  '

$ws.Range('D30').Value = '''1.32'
$ws.Range('E30').Value = '  -4.50%  '

$ws.Range('E31').Value = '  -4.59%  '

$ws.Range('E32').Value = '  -4.36%  '

$ws.Range('E33').Value = '  -4.75%  '

$ws.Range('E34').Value = '  -2.82%  '

$ws.Range('E35').Value = '  -4.89%  '

$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.12%  '

$ws.Range('D38').Value = '''3.17'
$ws.Range('E38').Value = '  -4.69%  '

$ws.Range('D39').Value = '''5.35'
$ws.Range('E39').Value = '  -0.37%  '

$ws.Range('D40').Value = '1.497.42'
$ws.Range('E40').Value = '  +3.12%  '

$ws.Range('E41').Value = '  -7.20%  '

$ws.Range('E42').Value = '  -2.07%  '

$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''95.06'
$ws.Range('E43').Value = '  -6.12%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0928'
$ws.Range('E44').Value = '  -4.00%  '

$ws.Range('D45').Value = '''16.56'
$ws.Range('E45').Value = '  -0.63%  '

$ws.Range('E46').Value = '  -5.28%  '

$ws.Range('E47').Value = '  -4.93%  '

$ws.Range('D48').Value = '''7.14'
$ws.Range('E48').Value = '  -0.92%  '

$ws.Range('E49').Value = '  -1.77%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.212.05'
$ws.Range('E50').Value = '  -2.87%  '

$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').Value = '''3.61'
$ws.Range('E51').Value = '  -12.13%  '
